$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'35.572.45"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.49%  "

# Row 3
$ws.Range("D3").Value = "'1.916.24"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.43%  "

# Row 4
$ws.Range("E4").Value = "  +0.56%  "

# Row 5
$ws.Range("D5").Value = "'246.78"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.11%  "

# Row 6
$ws.Range("D6").Value = "'0.655"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.22%  "

# Row 7
$ws.Range("E7").Value = "  +0.47%  "

# Row 8
$ws.Range("D8").Value = "'42.11"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.98%  "

# Row 9
$ws.Range("E9").Value = "  +5.99%  "

# Row 10
$ws.Range("D10").Value = "'49.44"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.84%  "

# Row 11
$ws.Range("D11").Value = "'0.0722"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.91%  "

# Row 12
$ws.Range("E12").Value = "  +1.12%  "

# Row 13
$ws.Range("D13").Value = "'2.196.91"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.53%  "

# Row 14
$ws.Range("D14").Value = "'12.32"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +7.66%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "'1.940.04"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.66%  "

# Row 16
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.702"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.45%  "

# Row 17
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'4.93"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.82%  "

# Row 18
$ws.Range("D18").Value = "'35.591.85"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.55%  "

# Row 19
$ws.Range("D19").Value = "'72.41"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.98%  "

# Row 20
$ws.Range("D20").Value = "'0.0₃0826"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.71%  "

# Row 21
$ws.Range("D21").Value = "'246.16"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.23%  "

# Row 22
$ws.Range("D22").Value = "'12.69"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.94%  "

# Row 23
$ws.Range("D23").Value = "'4.86"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.75%  "

# Row 24
$ws.Range("E24").Value = "  +0.50%  "

# Row 25
$ws.Range("D25").Value = "'2.31"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.09%  "

# Row 26
$ws.Range("D26").Value = "'2.21"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +16.54%  "

# Row 27
$ws.Range("D27").Value = "'171.71"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.55%  "

# Row 28
$ws.Range("D28").Value = "'8.44"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +5.65%  "

# Row 29
$ws.Range("D29").Value = "'18.61"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.20%  "

# Row 30
$ws.Range("E30").Value = "  +2.72%  "

# Row 31
$ws.Range("D31").Value = "'4.18"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.92%  "

# Row 32
$ws.Range("D32").Value = "'0.0573"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.64%  "

# Row 33
$ws.Range("B33").Value = "BinanceUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D33").Value = "'1.01"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.53%  "

# Row 34
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.936"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +19.28%  "

# Row 35
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "'4.18"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.01%  "

# Row 36
$ws.Range("D36").Value = "'1.76"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.21%  "

# Row 37
$ws.Range("D37").Value = "'2.05"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.60%  "

# Row 38
$ws.Range("E38").Value = "  +0.95%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0213"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.88%  "

# Row 40
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'1.11"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.89%  "

# Row 41
$ws.Range("E41").Value = "  +15.92%  "

# Row 42
$ws.Range("D42").Value = "'91.79"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.11%  "

# Row 43
$ws.Range("D43").Value = "'15.89"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +7.60%  "

# Row 44
$ws.Range("D44").Value = "'1.359.97"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.67%  "

# Row 45
$ws.Range("D45").Value = "'2.41"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.81%  "

# Row 46
$ws.Range("D46").Value = "'47.81"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +38.48%  "

# Row 47
$ws.Range("D47").Value = "'12.75"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.28%  "

# Row 48
$ws.Range("E48").Value = "  +2.44%  "

# Row 49
$ws.Range("D49").Value = "'2.41"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.17%  "

# Row 50
$ws.Range("D50").Value = "'6.59"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.94%  "

# Row 51
$ws.Range("D51").Value = "'2.107.28"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.51%  "
